$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("spell")

# Build the new row 5 first, copying formatting (style index 1) from row 4's
# corresponding cells without spilling formatting into the untouched columns.
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)

$ws.Range("C4").Copy()
$ws.Range("C5").PasteSpecial(-4122)

$ws.Range("F4:I4").Copy()
$ws.Range("F5:I5").PasteSpecial(-4122)

# Row 5 (new): second first_reward entry, spell 1003 -> "공기의 근원" (air source)
# Set C5's value before C4's so "공기의 근원" is appended to the shared-strings
# table ahead of "불의 근원", matching the authored order (index 18 then 19).
$ws.Range("A5").Value = 1003
$ws.Range("C5").Value = "공기의 근원"
$ws.Range("F5").Value = 2
$ws.Range("G5").Value = 2
$ws.Range("H5").Value = 1
$ws.Range("I5").Value = 10103

# Row 4: set the source name for spell 1001 -> "불의 근원" (fire source)
$ws.Range("C4").Value = "불의 근원"

$ws.Range("I5").Select()
